# Daily attendance processing - 2026-01-09 12:53:29
# Reorders the "Recorded By" (column G) names for a specific set of rows
# so that "System" is listed first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = "System, backup@backdoor.com, system"
    3   = "System, dnasr281@gmail.com"
    4   = "System, backup@backdoor.com"
    5   = "System, backup@backdoor.com"
    6   = "System, dnasr281@gmail.com"
    7   = "System, admin@admin.com"
    8   = "System, backup@backdoor.com"
    28  = "System, backup@backdoor.com, system"
    29  = "System, dnasr281@gmail.com"
    30  = "System, backup@backdoor.com"
    31  = "System, backup@backdoor.com"
    32  = "System, dnasr281@gmail.com"
    33  = "System, admin@admin.com"
    34  = "System, backup@backdoor.com"
    54  = "System, backup@backdoor.com, system"
    55  = "System, dnasr281@gmail.com"
    56  = "System, backup@backdoor.com"
    57  = "System, backup@backdoor.com"
    58  = "System, dnasr281@gmail.com"
    59  = "System, admin@admin.com"
    60  = "System, backup@backdoor.com"
    80  = "System, backup@backdoor.com"
    81  = "System, backup@backdoor.com"
    82  = "System, backup@backdoor.com"
    87  = "dnasr281@gmail.com, admin@admin.com"
    106 = "System, backup@backdoor.com"
    107 = "System, backup@backdoor.com"
    108 = "System, backup@backdoor.com"
    113 = "dnasr281@gmail.com, admin@admin.com"
    132 = "System, backup@backdoor.com"
    133 = "System, backup@backdoor.com"
    134 = "System, backup@backdoor.com"
    139 = "dnasr281@gmail.com, admin@admin.com"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
